$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (this also updates the Print_Area defined name
# that references the sheet by name).
$ws.Name = "IEEE talks 19"

# Select the frozen-pane worksheet window so we can adjust scroll/selection.
$win = $excel.ActiveWindow

# Scroll the bottom pane so row 16 becomes the first visible row below
# the frozen rows (previously row 4).
$win.ScrollRow = 16
$win.ScrollColumn = 1

# Update the active selection in the bottom-left pane to a single cell.
$ws.Range("L11").Select()
